$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.406999999999999
$ws.Range("B14").Value = 6.645999999999999
$ws.Range("B16").Value = 6.406999999999999
$ws.Range("B21").Value = 6.179
$ws.Range("B23").Value = 6.665000000000001
$ws.Range("B25").Value = 6.38
